# Apply the cryptos-list data refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that looks numeric ("228.82", "0.0621", ...).
# Pre-format the whole column as Text so the assignment below is not
# auto-coerced into a Number by the COM layer (matching the source file,
# where these are plain inline strings, not numeric cells).
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "39.476.76"
$ws.Range("E2").Value = "  +1.79%  "

$ws.Range("D3").Value = "2.165.96"
$ws.Range("E3").Value = "  +2.90%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "228.82"
$ws.Range("E5").Value = "  +0.44%  "

$ws.Range("D6").Value = "0.626"
$ws.Range("E6").Value = "  +1.42%  "

$ws.Range("D7").Value = "63.68"
$ws.Range("E7").Value = "  +2.34%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").Value = "0.396"
$ws.Range("E9").Value = "  +1.25%  "

$ws.Range("D10").Value = "0.0853"
$ws.Range("E10").Value = "  +1.34%  "

$ws.Range("E11").Value = "  +0.18%  "

$ws.Range("D12").Value = "16.08"
$ws.Range("E12").Value = "  +2.10%  "

$ws.Range("D13").Value = "2.486.38"
$ws.Range("E13").Value = "  +2.92%  "

$ws.Range("D14").Value = "22.04"
$ws.Range("E14").Value = "  -0.22%  "

$ws.Range("D15").Value = "0.813"
$ws.Range("E15").Value = "  +0.55%  "

$ws.Range("D16").Value = "5.53"
$ws.Range("E16").Value = "  +0.01%  "

$ws.Range("D17").Value = "2.152.55"
$ws.Range("E17").Value = "  +2.32%  "

$ws.Range("D18").Value = "39.479.32"
$ws.Range("E18").Value = "  +1.77%  "

$ws.Range("D19").Value = "6.24"
$ws.Range("E19").Value = "  +1.96%  "

$ws.Range("D20").Value = "71.90"
$ws.Range("E20").Value = "  -0.12%  "

$ws.Range("D21").Value = "0.0₃0848"
$ws.Range("E21").Value = "  +0.85%  "

$ws.Range("D22").Value = "229.63"
$ws.Range("E22").Value = "  +0.63%  "

$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").Value = "2.36"
$ws.Range("E24").Value = "  +1.65%  "

$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "2.34"
$ws.Range("E25").Value = "  -0.71%  "

$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "172.14"
$ws.Range("E26").Value = "  -0.13%  "

$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "9.56"
$ws.Range("E27").Value = "  -0.93%  "

$ws.Range("E28").Value = "  -0.69%  "

$ws.Range("D29").Value = "19.88"
$ws.Range("E29").Value = "  +2.70%  "

$ws.Range("D30").Value = "1.42"
$ws.Range("E30").Value = "  +0.10%  "

$ws.Range("D31").Value = "2.64"
$ws.Range("E31").Value = "  +4.66%  "

$ws.Range("E32").Value = "  +1.23%  "

$ws.Range("D33").Value = "4.62"
$ws.Range("E33").Value = "  +1.47%  "

$ws.Range("E34").Value = "  -1.05%  "

$ws.Range("D35").Value = "7.06"
$ws.Range("E35").Value = "  +0.29%  "

$ws.Range("D36").Value = "0.0621"
$ws.Range("E36").Value = "  +0.02%  "

$ws.Range("D37").Value = "2.45"
$ws.Range("E37").Value = "  +1.53%  "

$ws.Range("D38").Value = "3.62"
$ws.Range("E38").Value = "  +0.69%  "

$ws.Range("E39").Value = "  -0.15%  "

$ws.Range("D40").Value = "103.09"
$ws.Range("E40").Value = "  +0.33%  "

$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").Value = "17.97"
$ws.Range("E41").Value = "  -0.60%  "

$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "0.0228"
$ws.Range("E42").Value = "  +0.05%  "

$ws.Range("D43").Value = "1.525.16"
$ws.Range("E43").Value = "  -0.25%  "

$ws.Range("E44").Value = "  +0.41%  "

$ws.Range("E45").Value = "  +5.26%  "

$ws.Range("E46").Value = "  +0.64%  "

$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "0.0925"
$ws.Range("E47").Value = "  +1.41%  "

$ws.Range("B48").Value = "FTXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D48").Value = "4.27"
$ws.Range("E48").Value = "  +3.49%  "

$ws.Range("D49").Value = "7.77"
$ws.Range("E49").Value = "  +0.16%  "

$ws.Range("D50").Value = "2.369.43"
$ws.Range("E50").Value = "  +2.93%  "

$ws.Range("B51").Value = "SynthetixNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D51").Value = "3.89"
$ws.Range("E51").Value = "  +9.38%  "

# Restore the default (General/no explicit style) formatting on column D
# now that the text values are locked in, so cells keep their original
# unstyled appearance.
$dRange.Style = "Normal"